# Realestate Update resale numbers 2023-06-10 17:04
# Append a new data row (row 38) with the latest resale numbers snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

# Date/Week columns are stored as plain text in this sheet (not real
# Excel dates/numbers). Temporarily force a text number-format before
# assigning the value so Excel's auto-detection doesn't turn them into
# a date serial / numeric value, then clear the format again so the
# cell doesn't end up with a stray style index.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-10"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "16:58:53"
$ws.Cells.Item($row, 3).Value = "Saturday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "23"
$ws.Cells.Item($row, 4).ClearFormats()

# Per-city resale numbers (numeric columns)
$ws.Cells.Item($row, 5).Value = 121176
$ws.Cells.Item($row, 6).Value = 134626
$ws.Cells.Item($row, 7).Value = 160934
$ws.Cells.Item($row, 8).Value = 132019
$ws.Cells.Item($row, 9).Value = 176133
$ws.Cells.Item($row, 10).Value = 114025
$ws.Cells.Item($row, 11).Value = 201893
$ws.Cells.Item($row, 12).Value = 222228
$ws.Cells.Item($row, 13).Value = 173597
$ws.Cells.Item($row, 14).Value = 98513
$ws.Cells.Item($row, 15).Value = 38730
$ws.Cells.Item($row, 16).Value = 34283
$ws.Cells.Item($row, 17).Value = 51156
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36948
$ws.Cells.Item($row, 20).Value = -1
